$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Avverkningsanmälningar")

# Update the "Förändrad" date column (C2:C27) from 2024-01-21 to 2024-01-24
for ($r = 2; $r -le 27; $r++) {
    $ws.Cells.Item($r, 3).Value = 45315
}

# Delete the last data row (row 28: "A 2414-2024")
$ws.Rows.Item(28).Delete()
